# Generate Report for Handback
# Regenerates the handback-status report: new source-file GUIDs
# (981f8289-...  -> c3a5706a-...  and dd4315f8-... -> fffff01fc3cf-...),
# a new handoff/handback xlf hash (dc122125...-> fe5eebdf8a...), and
# refreshed handoff/handback timestamps, across the Overview / zh-cn / de-de
# sheets. Hyperlink targets (the github blob URLs) are unchanged.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": only the two hyperlinked file-name cells change.
# ---------------------------------------------------------------------
$ws1 = $wb.Sheets.Item("Overview")
$ws1.Range("A1").Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/9651d68e8c0dcb9ad3985bd8d4df2a0a84929d2b/e2e/981f8289-ba26-4bb1-bd0f-9ee34b00cffb.md", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/9651d68e8c0dcb9ad3985bd8d4df2a0a84929d2b/e2e/dd4315f8-9c62-41a2-af39-71724def01b7.md", `
    "", "", "fffff01fc3cf-055b-49b3-bde0-c5320a947af2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Sheets.Item("zh-cn")
$ws2.Range("A1").Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/9651d68e8c0dcb9ad3985bd8d4df2a0a84929d2b/e2e/981f8289-ba26-4bb1-bd0f-9ee34b00cffb.md", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ee759d7418e0671790d0a288faf94677bf07b5ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/981f8289-ba26-4bb1-bd0f-9ee34b00cffb.dc122125df8f338292684c4a98312cb5ac70ecc7.zh-cn.xlf", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.fe5eebdf8a1eac2b97594558ab8ae40a84c285b1.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), `
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/73a22ad96ffeacc492268b7beea59485bf94952a/e2e/981f8289-ba26-4bb1-bd0f-9ee34b00cffb.md", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/880ddac4ebcb9830c484460ae96a74dc0ef04bcf/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/981f8289-ba26-4bb1-bd0f-9ee34b00cffb.dc122125df8f338292684c4a98312cb5ac70ecc7.zh-cn.xlf", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.fe5eebdf8a1eac2b97594558ab8ae40a84c285b1.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/9651d68e8c0dcb9ad3985bd8d4df2a0a84929d2b/e2e/dd4315f8-9c62-41a2-af39-71724def01b7.md", `
    "", "", "fffff01fc3cf-055b-49b3-bde0-c5320a947af2.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ee759d7418e0671790d0a288faf94677bf07b5ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/dd4315f8-9c62-41a2-af39-71724def01b7.d87f492af258b291bfefcfdbf39da92b0792710b.zh-cn.xlf", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.fe5eebdf8a1eac2b97594558ab8ae40a84c285b1.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), `
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/73a22ad96ffeacc492268b7beea59485bf94952a/e2e/dd4315f8-9c62-41a2-af39-71724def01b7.md", `
    "", "", "fffff01fc3cf-055b-49b3-bde0-c5320a947af2.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/880ddac4ebcb9830c484460ae96a74dc0ef04bcf/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/dd4315f8-9c62-41a2-af39-71724def01b7.d87f492af258b291bfefcfdbf39da92b0792710b.zh-cn.xlf", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.fe5eebdf8a1eac2b97594558ab8ae40a84c285b1.zh-cn.xlf")

# Handoff / handback timestamps (plain text cells, no hyperlink)
$ws2.Range("E2").Value = "2016-03-21 23:01:39"
$ws2.Range("H2").Value = "2016-03-21 23:02:09"
$ws2.Range("E3").Value = "2016-03-21 23:01:39"
$ws2.Range("H3").Value = "2016-03-21 23:02:09"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Sheets.Item("de-de")
$ws3.Range("A1").Hyperlinks.Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/9651d68e8c0dcb9ad3985bd8d4df2a0a84929d2b/e2e/981f8289-ba26-4bb1-bd0f-9ee34b00cffb.md", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf57fbd3773a8cb059c94af650b9881e52ba7b49/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/981f8289-ba26-4bb1-bd0f-9ee34b00cffb.dc122125df8f338292684c4a98312cb5ac70ecc7.de-de.xlf", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.fe5eebdf8a1eac2b97594558ab8ae40a84c285b1.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), `
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8e13afaad00e732441c6b03ed74022e3a3430405/e2e/981f8289-ba26-4bb1-bd0f-9ee34b00cffb.md", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b524fbd6591a19782c96fd9416589f1cfb9b006c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/981f8289-ba26-4bb1-bd0f-9ee34b00cffb.dc122125df8f338292684c4a98312cb5ac70ecc7.de-de.xlf", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.fe5eebdf8a1eac2b97594558ab8ae40a84c285b1.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), `
    "https://github.com/OpenLocalizationTest/oltest/blob/9651d68e8c0dcb9ad3985bd8d4df2a0a84929d2b/e2e/dd4315f8-9c62-41a2-af39-71724def01b7.md", `
    "", "", "fffff01fc3cf-055b-49b3-bde0-c5320a947af2.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf57fbd3773a8cb059c94af650b9881e52ba7b49/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/dd4315f8-9c62-41a2-af39-71724def01b7.d87f492af258b291bfefcfdbf39da92b0792710b.de-de.xlf", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.fe5eebdf8a1eac2b97594558ab8ae40a84c285b1.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), `
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8e13afaad00e732441c6b03ed74022e3a3430405/e2e/dd4315f8-9c62-41a2-af39-71724def01b7.md", `
    "", "", "fffff01fc3cf-055b-49b3-bde0-c5320a947af2.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b524fbd6591a19782c96fd9416589f1cfb9b006c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/dd4315f8-9c62-41a2-af39-71724def01b7.d87f492af258b291bfefcfdbf39da92b0792710b.de-de.xlf", `
    "", "", "c3a5706a-ec37-4ce4-899a-1de5f0dfe347.fe5eebdf8a1eac2b97594558ab8ae40a84c285b1.de-de.xlf")

# Handoff / handback timestamps (plain text cells, no hyperlink)
$ws3.Range("E2").Value = "2016-03-21 23:01:44"
$ws3.Range("H2").Value = "2016-03-21 23:02:15"
$ws3.Range("E3").Value = "2016-03-21 23:01:44"
$ws3.Range("H3").Value = "2016-03-21 23:02:15"

Write-Host "Handback report regenerated."
